$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous periods (2204-2210, rows 16-22) are replaced with the new
# periods sorted in descending order (2210 down to 2204), carrying their
# "Valor Mora" amounts along with them.
$periods = @("2210", "2209", "2208", "2207", "2206", "2205", "2204")
$valores = @(32707, 40000, 40000, 40000, 40000, 40000, 40000)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
